$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2457.2173
$ws.Range("I2").Value = 3247.9375
$ws.Range("J2").Value = 649.8570999999999
$ws.Range("K2").Value = 3247.9375
$ws.Range("L2").Value = 649.8570999999999
$ws.Range("M2").Value = -3134.9375
$ws.Range("N2").Value = -875.8570999999999
$ws.Range("H33").Value = 235.9
$ws.Range("I33").Value = 249.5
$ws.Range("K33").Value = 249.5
$ws.Range("M33").Value = -20.5
$ws.Range("H40").Value = 4959.7
$ws.Range("J40").Value = 6458.625
$ws.Range("L40").Value = 6458.625
$ws.Range("N40").Value = -6808.625
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H116").Value = 5484.5835
$ws.Range("I116").Value = 5103.1113
$ws.Range("K116").Value = 5103.1113
$ws.Range("M116").Value = -1661.1113
$ws.Range("H121").Value = 981.5
$ws.Range("J121").Value = 981.5
$ws.Range("L121").Value = 2944.5
$ws.Range("N121").Value = -6438.5
$ws.Range("H125").Value = 2499.5
$ws.Range("I125").Value = 3000
$ws.Range("J125").Value = 1999
$ws.Range("K125").Value = 27000
$ws.Range("L125").Value = 17991
$ws.Range("M125").Value = -24540
$ws.Range("N125").Value = -22911
$ws.Range("H131").Value = 1097.8
$ws.Range("I131").Value = 1097.8
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 3293.4
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 1746.6
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2059.3572
$ws.Range("I132").Value = 1209.7407
$ws.Range("K132").Value = 3629.2221
$ws.Range("M132").Value = -1099.2221
$ws.Range("H137").Value = 3250.95
$ws.Range("I137").Value = 1668.6
$ws.Range("K137").Value = 5005.799999999999
$ws.Range("M137").Value = -2455.799999999999
$ws.Range("H138").Value = 2245.8462
$ws.Range("I138").Value = 1409.3334
$ws.Range("K138").Value = 4228.0002
$ws.Range("M138").Value = 911.9997999999996
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2772
$ws.Range("I2").Value = 2772
$ws.Range("K2").Value = 2772
$ws.Range("M2").Value = -2659
$ws.Range("H32").Value = 24137.637
$ws.Range("I32").Value = 25500
$ws.Range("J32").Value = 23002.334
$ws.Range("K32").Value = 25500
$ws.Range("L32").Value = 23002.334
$ws.Range("M32").Value = -25213
$ws.Range("N32").Value = -23576.334
$ws.Range("H63").Value = 4726.5
$ws.Range("I63").Value = 2627.4285
$ws.Range("K63").Value = 2627.4285
$ws.Range("M63").Value = -1941.4285
$ws.Range("H66").Value = 4726.5
$ws.Range("I66").Value = 2627.4285
$ws.Range("K66").Value = 13137.1425
$ws.Range("M66").Value = -9705.1425
$ws.Range("H74").Value = 2038.875
$ws.Range("I74").Value = 1839.862
$ws.Range("K74").Value = 1839.862
$ws.Range("M74").Value = -965.8620000000001
$ws.Range("H77").Value = 2038.875
$ws.Range("I77").Value = 1839.862
$ws.Range("K77").Value = 9199.310000000001
$ws.Range("M77").Value = -4831.310000000001
$ws.Range("H116").Value = 2772
$ws.Range("I116").Value = 2772
$ws.Range("K116").Value = 2772
$ws.Range("M116").Value = -478
$ws.Range("H122").Value = 4245.25
$ws.Range("I122").Value = 4245.25
$ws.Range("K122").Value = 12735.75
$ws.Range("M122").Value = -10285.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 58963
$ws.Range("J2").Value = 58963
$ws.Range("L2").Value = 58963
$ws.Range("N2").Value = -59189
$ws.Range("H3").Value = 2772
$ws.Range("I3").Value = 2772
$ws.Range("K3").Value = 2772
$ws.Range("M3").Value = -2658
$ws.Range("H7").Value = 12500719
$ws.Range("I7").Value = 9444786
$ws.Range("K7").Value = 9444786
$ws.Range("M7").Value = -9444673
$ws.Range("H13").Value = 58966
$ws.Range("J13").Value = 58966
$ws.Range("L13").Value = 58966
$ws.Range("N13").Value = -59302
$ws.Range("H20").Value = 6845
$ws.Range("I20").Value = 8736
$ws.Range("K20").Value = 8736
$ws.Range("M20").Value = -8489
$ws.Range("H64").Value = 583.5
$ws.Range("I64").Value = 318
$ws.Range("K64").Value = 318
$ws.Range("M64").Value = -93
$ws.Range("H67").Value = 583.5
$ws.Range("I67").Value = 318
$ws.Range("K67").Value = 318
$ws.Range("M67").Value = 462
$ws.Range("H80").Value = 196.57895
$ws.Range("I80").Value = 173.16667
$ws.Range("K80").Value = 173.16667
$ws.Range("M80").Value = 824.8333299999999
$ws.Range("H83").Value = 196.57895
$ws.Range("I83").Value = 173.16667
$ws.Range("K83").Value = 865.8333500000001
$ws.Range("M83").Value = 4126.16665
$ws.Range("H86").Value = 8150
$ws.Range("I86").Value = 7200
$ws.Range("J86").Value = 9100
$ws.Range("K86").Value = 7200
$ws.Range("L86").Value = 9100
$ws.Range("M86").Value = -6077
$ws.Range("N86").Value = -11346
$ws.Range("H89").Value = 8150
$ws.Range("I89").Value = 7200
$ws.Range("J89").Value = 9100
$ws.Range("K89").Value = 36000
$ws.Range("L89").Value = 45500
$ws.Range("M89").Value = -30384
$ws.Range("N89").Value = -56732
$ws.Range("H105").Value = 2759.9092
$ws.Range("I105").Value = 2181.8
$ws.Range("J105").Value = 3241.6667
$ws.Range("K105").Value = 2181.8
$ws.Range("L105").Value = 3241.6667
$ws.Range("M105").Value = -434.8000000000002
$ws.Range("N105").Value = -6735.6667
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6624.712
$ws.Range("I31").Value = 4137.6577
$ws.Range("K31").Value = 4137.6577
$ws.Range("M31").Value = -3842.6577
$ws.Range("H34").Value = 6624.712
$ws.Range("I34").Value = 4137.6577
$ws.Range("K34").Value = 4137.6577
$ws.Range("M34").Value = -3935.6577
$ws.Range("H99").Value = 5603.8887
$ws.Range("I99").Value = 5063.7144
$ws.Range("K99").Value = 5063.7144
$ws.Range("M99").Value = -3565.7144
$ws.Range("H107").Value = 422.9091
$ws.Range("I107").Value = 355.27777
$ws.Range("J107").Value = 727.25
$ws.Range("K107").Value = 355.27777
$ws.Range("L107").Value = 727.25
$ws.Range("M107").Value = 1564.72223
$ws.Range("N107").Value = -4567.25
$ws.Range("H126").Value = 5603.8887
$ws.Range("I126").Value = 5063.7144
$ws.Range("K126").Value = 15191.1432
$ws.Range("M126").Value = -12721.1432
$ws.Range("H134").Value = 2074.4866
$ws.Range("I134").Value = 1854.0883
$ws.Range("K134").Value = 5562.2649
$ws.Range("M134").Value = -3027.2649
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 134.2
$ws.Range("I23").Value = 130
$ws.Range("K23").Value = 390
$ws.Range("M23").Value = -155
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3499.5
$ws.Range("I126").Value = 3499.5
$ws.Range("K126").Value = 10498.5
$ws.Range("M126").Value = -8028.5
$ws.Range("H132").Value = 88861.30499999999
$ws.Range("I132").Value = 102745.37
$ws.Range("K132").Value = 308236.11
$ws.Range("M132").Value = -305706.11
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6245.6
$ws.Range("I7").Value = 5567.2
$ws.Range("K7").Value = 5567.2
$ws.Range("M7").Value = -5455.2
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H55").Value = 1891.4117
$ws.Range("J55").Value = 1907.9166
$ws.Range("L55").Value = 1907.9166
$ws.Range("N55").Value = -2253.9166
$ws.Range("H68").Value = 9399.4
$ws.Range("J68").Value = 9666.666999999999
$ws.Range("L68").Value = 9666.666999999999
$ws.Range("N68").Value = -11164.667
$ws.Range("H71").Value = 9399.4
$ws.Range("J71").Value = 9666.666999999999
$ws.Range("L71").Value = 48333.335
$ws.Range("N71").Value = -55821.335
$ws.Range("H126").Value = 6245.6
$ws.Range("I126").Value = 5567.2
$ws.Range("K126").Value = 16701.6
$ws.Range("M126").Value = -14231.6
$ws.Range("H132").Value = 8726.823
$ws.Range("I132").Value = 6948.5
$ws.Range("K132").Value = 20845.5
$ws.Range("M132").Value = -18315.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 7416.6665
$ws.Range("I61").Value = 7416.6665
$ws.Range("K61").Value = 7416.6665
$ws.Range("M61").Value = -7124.6665
$ws.Range("H96").Value = 1833
$ws.Range("J96").Value = 2640
$ws.Range("L96").Value = 2640
$ws.Range("N96").Value = -5386
$ws.Range("H113").Value = 833.1111
$ws.Range("J113").Value = 819.6
$ws.Range("L113").Value = 2458.8
$ws.Range("N113").Value = -6798.8
$ws.Range("H126").Value = 7212.5
$ws.Range("I126").Value = 4925
$ws.Range("K126").Value = 14775
$ws.Range("M126").Value = -12305

Write-Output "Applied all cell updates"